$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to store the value as text, matching the inline/shared string
    # representation in the original workbook (avoids Excel auto-converting
    # numeric-looking strings like "1.00" or "26.293.65" into numbers).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.293.65"
Set-TextValue $ws.Range("D3") "1.619.36"
Set-TextValue $ws.Range("E3") "  +2.00%  "
Set-TextValue $ws.Range("E4") "  -0.08%  "
Set-TextValue $ws.Range("D5") "211.94"
Set-TextValue $ws.Range("E5") "  +0.74%  "
Set-TextValue $ws.Range("E6") "  -0.09%  "
Set-TextValue $ws.Range("D7") "0.483"
Set-TextValue $ws.Range("E7") "  +0.95%  "
Set-TextValue $ws.Range("E8") "  +0.73%  "
Set-TextValue $ws.Range("D9") "0.0615"
Set-TextValue $ws.Range("E9") "  +0.81%  "
Set-TextValue $ws.Range("D10") "18.79"
Set-TextValue $ws.Range("E10") "  +4.88%  "
Set-TextValue $ws.Range("E11") "  +1.04%  "
Set-TextValue $ws.Range("D12") "1.844.40"
Set-TextValue $ws.Range("D13") "1.608.86"
Set-TextValue $ws.Range("E13") "  +1.35%  "
Set-TextValue $ws.Range("E14") "  +0.71%  "
Set-TextValue $ws.Range("D15") "0.518"
Set-TextValue $ws.Range("E15") "  +1.54%  "
Set-TextValue $ws.Range("D16") "26.300.97"
Set-TextValue $ws.Range("E16") "  +1.34%  "
Set-TextValue $ws.Range("D17") "62.24"
Set-TextValue $ws.Range("E17") "  +3.73%  "
Set-TextValue $ws.Range("D18") "0.0₃0727"
Set-TextValue $ws.Range("E18") "  +1.01%  "
Set-TextValue $ws.Range("E19") "  -0.07%  "
Set-TextValue $ws.Range("D20") "201.61"
Set-TextValue $ws.Range("E20") "  +1.32%  "
Set-TextValue $ws.Range("E21") "  +1.79%  "
Set-TextValue $ws.Range("D22") "9.32"
Set-TextValue $ws.Range("E22") "  +1.66%  "
Set-TextValue $ws.Range("D23") "6.05"
Set-TextValue $ws.Range("E23") "  +1.36%  "
Set-TextValue $ws.Range("E24") "  +3.27%  "
Set-TextValue $ws.Range("D25") "144.53"
Set-TextValue $ws.Range("E25") "  +1.37%  "
Set-TextValue $ws.Range("E26") "  -0.05%  "
Set-TextValue $ws.Range("E27") "  -1.33%  "
Set-TextValue $ws.Range("E28") "  +0.80%  "
Set-TextValue $ws.Range("D29") "6.55"
Set-TextValue $ws.Range("E29") "  +1.72%  "
Set-TextValue $ws.Range("D30") "0.0520"
Set-TextValue $ws.Range("E30") "  +10.18%  "
Set-TextValue $ws.Range("E31") "  +1.02%  "
Set-TextValue $ws.Range("E32") "  +1.97%  "
Set-TextValue $ws.Range("D33") "2.93"
Set-TextValue $ws.Range("E33") "  +0.03%  "
Set-TextValue $ws.Range("D34") "1.50"
Set-TextValue $ws.Range("E34") "  +2.07%  "
Set-TextValue $ws.Range("E35") "  +1.80%  "
Set-TextValue $ws.Range("D36") "1.177.60"
Set-TextValue $ws.Range("E36") "  +4.96%  "
Set-TextValue $ws.Range("D37") "0.0164"
Set-TextValue $ws.Range("E37") "  +0.71%  "
Set-TextValue $ws.Range("D38") "0.806"
Set-TextValue $ws.Range("E38") "  +3.15%  "
Set-TextValue $ws.Range("E39") "  -0.07%  "
Set-TextValue $ws.Range("E40") "  +0.21%  "
Set-TextValue $ws.Range("D41") "0.495"
Set-TextValue $ws.Range("E41") "  +1.70%  "
Set-TextValue $ws.Range("E42") "  +1.69%  "
Set-TextValue $ws.Range("D43") "5.34"
Set-TextValue $ws.Range("E43") "  +5.04%  "
Set-TextValue $ws.Range("D44") "1.755.19"
Set-TextValue $ws.Range("E44") "  +2.02%  "
Set-TextValue $ws.Range("D45") "92.70"
Set-TextValue $ws.Range("E45") "  +1.00%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "1.54"
Set-TextValue $ws.Range("E46") "  +3.85%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "53.74"
Set-TextValue $ws.Range("E47") "  +1.08%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.0508"
Set-TextValue $ws.Range("E48") "  +1.08%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.408"
Set-TextValue $ws.Range("E49") "  +0.40%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue $ws.Range("D50") "1.00"
Set-TextValue $ws.Range("E50") "  -0.14%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.26"
Set-TextValue $ws.Range("E51") "  +2.59%  "
